$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Rows 78 and 79 had their match data swapped (the "Betis vs
#    Valencia" match and the "Atl. Madrid vs Cadiz CF" match traded
#    places). Columns A-E (index, pais, torneio, temporada, data) and
#    G (home_ft_gols) stay identical between the two rows, only
#    F, H, I, J-V differ, so overwrite those explicitly.
# -----------------------------------------------------------------

# New content for row 78 (previously the Betis/Valencia match, now
# becomes the Atl. Madrid / Cadiz CF match)
$ws.Range("F78").Value = "Atl. Madrid"
$ws.Range("H78").Value = "Cadiz CF"
$ws.Range("I78").Value = 2
$ws.Range("J78").Value = 1.35
$ws.Range("K78").Value = "21/09/2023 22:03"
$ws.Range("L78").Value = 1.34
$ws.Range("M78").Value = "01/10/2023 20:50"
$ws.Range("N78").Value = 5.24
$ws.Range("O78").Value = "21/09/2023 22:03"
$ws.Range("P78").Value = 5.31
$ws.Range("Q78").Value = "01/10/2023 20:59"
$ws.Range("R78").Value = 9.25
$ws.Range("S78").Value = "21/09/2023 22:03"
$ws.Range("T78").Value = 10.48
$ws.Range("U78").Value = "01/10/2023 20:59"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/spain/laliga/atl-madrid-cadiz/E1cOKVAj/"

# New content for row 79 (previously the Atl. Madrid/Cadiz CF match,
# now becomes the Betis / Valencia match)
$ws.Range("F79").Value = "Betis"
$ws.Range("H79").Value = "Valencia"
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 2.5
$ws.Range("K79").Value = "24/09/2023 17:02"
$ws.Range("L79").Value = 2.07
$ws.Range("M79").Value = "01/10/2023 20:54"
$ws.Range("N79").Value = 3.19
$ws.Range("O79").Value = "24/09/2023 17:02"
$ws.Range("P79").Value = 3.41
$ws.Range("Q79").Value = "01/10/2023 20:54"
$ws.Range("R79").Value = 3.11
$ws.Range("S79").Value = "24/09/2023 17:02"
$ws.Range("T79").Value = 4.04
$ws.Range("U79").Value = "01/10/2023 20:57"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/spain/laliga/betis-valencia/vukArZ2c/"

# -----------------------------------------------------------------
# 2) A new row (122) with a new match (Rayo Vallecano vs Girona) is
#    appended after the previous last row (121). Copy the formatting
#    of row 121 down to row 122 first (so the index column keeps its
#    bold/bordered style and the date column keeps its date number
#    format), then fill in the values.
# -----------------------------------------------------------------
$ws.Range("A121:V121").Copy()
$ws.Range("A122:V122").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A122").Value = 121
$ws.Range("B122").Value = "spain"
$ws.Range("C122").Value = "laliga"
$ws.Range("D122").Value = "2023-2024"
$ws.Range("E122").Value = 45241.58333333334
$ws.Range("F122").Value = "Rayo Vallecano"
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = "Girona"
$ws.Range("I122").Value = 2
$ws.Range("J122").Value = 2.47
$ws.Range("K122").Value = "29/10/2023 11:02"
$ws.Range("L122").Value = 2.99
$ws.Range("M122").Value = "11/11/2023 13:58"
$ws.Range("N122").Value = 3.32
$ws.Range("O122").Value = "29/10/2023 11:02"
$ws.Range("P122").Value = 3.5
$ws.Range("Q122").Value = "11/11/2023 13:27"
$ws.Range("R122").Value = 2.99
$ws.Range("S122").Value = "29/10/2023 11:02"
$ws.Range("T122").Value = 2.46
$ws.Range("U122").Value = "11/11/2023 13:58"
$ws.Range("V122").Value = "https://www.betexplorer.com/football/spain/laliga/rayo-vallecano-girona/z3h7PEcB/"
